# Auto-generated edit script: updates COVID-19 country data table
# Source: diff between before.xlsx and after.xlsx (paises.xlsx / sheet 'Pais')
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(1, 1).Value = 'Datos actualizados a 29 de Marzo de 2020 a las 17:50'

$ws.Cells.Item(4, 1).Value = 'Estados Unidos'
$ws.Cells.Item(4, 2).Value = 125099
$ws.Cells.Item(4, 3).Value = 1521
$ws.Cells.Item(4, 4).Value = 3238
$ws.Cells.Item(4, 5).Value = 119623
$ws.Cells.Item(4, 6).Value = 2666
$ws.Cells.Item(4, 7).Value = 17
$ws.Cells.Item(4, 8).Value = 2238

$ws.Cells.Item(7, 1).Value = 'España'
$ws.Cells.Item(7, 2).Value = 78799
$ws.Cells.Item(7, 3).Value = 5564
$ws.Cells.Item(7, 4).Value = 14709
$ws.Cells.Item(7, 5).Value = 57484
$ws.Cells.Item(7, 6).Value = 4165
$ws.Cells.Item(7, 7).Value = 624
$ws.Cells.Item(7, 8).Value = 6606

$ws.Cells.Item(26, 1).Value = 'Malasia'
$ws.Cells.Item(26, 2).Value = 2470
$ws.Cells.Item(26, 3).Value = 150
$ws.Cells.Item(26, 4).Value = 388
$ws.Cells.Item(26, 5).Value = 2077
$ws.Cells.Item(26, 6).Value = 73
$ws.Cells.Item(26, 7).Value = 22
$ws.Cells.Item(26, 8).Value = 5

$ws.Cells.Item(30, 1).Value = 'Luxemburgo'
$ws.Cells.Item(30, 2).Value = 1950
$ws.Cells.Item(30, 3).Value = 119
$ws.Cells.Item(30, 4).Value = 40
$ws.Cells.Item(30, 5).Value = 1889
$ws.Cells.Item(30, 6).Value = 25
$ws.Cells.Item(30, 7).Value = 3
$ws.Cells.Item(30, 8).Value = 21

$ws.Cells.Item(43, 1).Value = 'Grecia'
$ws.Cells.Item(43, 2).Value = 1156
$ws.Cells.Item(43, 3).Value = 95
$ws.Cells.Item(43, 4).Value = 52
$ws.Cells.Item(43, 5).Value = 1066
$ws.Cells.Item(43, 6).Value = 66
$ws.Cells.Item(43, 7).Value = 6
$ws.Cells.Item(43, 8).Value = 38

$ws.Cells.Item(106, 1).Value = 'Estado de Palestina'
$ws.Cells.Item(106, 2).Value = 108
$ws.Cells.Item(106, 3).Value = 4
$ws.Cells.Item(106, 4).Value = 18
$ws.Cells.Item(106, 5).Value = 89
$ws.Cells.Item(106, 6).Value = 0
$ws.Cells.Item(106, 7).Value = 0
$ws.Cells.Item(106, 8).Value = 1

$ws.Cells.Item(107, 1).Value = 'Mauricio'
$ws.Cells.Item(107, 2).Value = 107
$ws.Cells.Item(107, 3).Value = 5
$ws.Cells.Item(107, 4).Value = 0
$ws.Cells.Item(107, 5).Value = 105
$ws.Cells.Item(107, 6).Value = 1
$ws.Cells.Item(107, 7).Value = 0
$ws.Cells.Item(107, 8).Value = 2

$ws.Cells.Item(114, 1).Value = 'Georgia'
$ws.Cells.Item(114, 2).Value = 91
$ws.Cells.Item(114, 3).Value = 1
$ws.Cells.Item(114, 4).Value = 18
$ws.Cells.Item(114, 5).Value = 73
$ws.Cells.Item(114, 6).Value = 1
$ws.Cells.Item(114, 7).Value = 0
$ws.Cells.Item(114, 8).Value = 0

$ws.Cells.Item(117, 1).Value = 'Bolivia'
$ws.Cells.Item(117, 2).Value = 81
$ws.Cells.Item(117, 3).Value = 7
$ws.Cells.Item(117, 4).Value = 0
$ws.Cells.Item(117, 5).Value = 80
$ws.Cells.Item(117, 6).Value = 3
$ws.Cells.Item(117, 7).Value = 1
$ws.Cells.Item(117, 8).Value = 1

$ws.Cells.Item(128, 1).Value = 'Isla de Man'
$ws.Cells.Item(128, 2).Value = 42
$ws.Cells.Item(128, 3).Value = 10
$ws.Cells.Item(128, 4).Value = 0
$ws.Cells.Item(128, 5).Value = 42
$ws.Cells.Item(128, 6).Value = 0
$ws.Cells.Item(128, 7).Value = 0
$ws.Cells.Item(128, 8).Value = 0

$ws.Cells.Item(129, 1).Value = 'Kenia'
$ws.Cells.Item(129, 2).Value = 42
$ws.Cells.Item(129, 3).Value = 4
$ws.Cells.Item(129, 4).Value = 1
$ws.Cells.Item(129, 5).Value = 40
$ws.Cells.Item(129, 6).Value = 2
$ws.Cells.Item(129, 7).Value = 0
$ws.Cells.Item(129, 8).Value = 1

$ws.Cells.Item(130, 1).Value = 'Madagascar'
$ws.Cells.Item(130, 2).Value = 39
$ws.Cells.Item(130, 3).Value = 13
$ws.Cells.Item(130, 4).Value = 0
$ws.Cells.Item(130, 5).Value = 39
$ws.Cells.Item(130, 6).Value = 0
$ws.Cells.Item(130, 7).Value = 0
$ws.Cells.Item(130, 8).Value = 0

$ws.Cells.Item(131, 1).Value = 'Puerto Rico'
$ws.Cells.Item(131, 2).Value = 39
$ws.Cells.Item(131, 3).Value = 0
$ws.Cells.Item(131, 4).Value = 1
$ws.Cells.Item(131, 5).Value = 36
$ws.Cells.Item(131, 6).Value = 0
$ws.Cells.Item(131, 7).Value = 0
$ws.Cells.Item(131, 8).Value = 2

$ws.Cells.Item(169, 1).Value = 'Libia'
$ws.Cells.Item(169, 2).Value = 8
$ws.Cells.Item(169, 3).Value = 5
$ws.Cells.Item(169, 4).Value = 0
$ws.Cells.Item(169, 5).Value = 8
$ws.Cells.Item(169, 6).Value = 0
$ws.Cells.Item(169, 7).Value = 0
$ws.Cells.Item(169, 8).Value = 0

$ws.Cells.Item(170, 1).Value = 'Guyana'
$ws.Cells.Item(170, 2).Value = 8
$ws.Cells.Item(170, 3).Value = 0
$ws.Cells.Item(170, 4).Value = 0
$ws.Cells.Item(170, 5).Value = 7
$ws.Cells.Item(170, 6).Value = 0
$ws.Cells.Item(170, 7).Value = 0
$ws.Cells.Item(170, 8).Value = 1

$ws.Cells.Item(171, 1).Value = 'Islas Caimanes'
$ws.Cells.Item(171, 2).Value = 8
$ws.Cells.Item(171, 3).Value = 0
$ws.Cells.Item(171, 4).Value = 0
$ws.Cells.Item(171, 5).Value = 7
$ws.Cells.Item(171, 6).Value = 0
$ws.Cells.Item(171, 7).Value = 0
$ws.Cells.Item(171, 8).Value = 1

$ws.Cells.Item(172, 1).Value = 'Curazao'
$ws.Cells.Item(172, 2).Value = 8
$ws.Cells.Item(172, 3).Value = 0
$ws.Cells.Item(172, 4).Value = 2
$ws.Cells.Item(172, 5).Value = 5
$ws.Cells.Item(172, 6).Value = 0
$ws.Cells.Item(172, 7).Value = 0
$ws.Cells.Item(172, 8).Value = 1

$ws.Cells.Item(173, 1).Value = 'Antigua y Barbuda'
$ws.Cells.Item(173, 2).Value = 7
$ws.Cells.Item(173, 3).Value = 0
$ws.Cells.Item(173, 4).Value = 0
$ws.Cells.Item(173, 5).Value = 7
$ws.Cells.Item(173, 6).Value = 0
$ws.Cells.Item(173, 7).Value = 0
$ws.Cells.Item(173, 8).Value = 0

$ws.Cells.Item(174, 1).Value = 'Gabon'
$ws.Cells.Item(174, 2).Value = 7
$ws.Cells.Item(174, 3).Value = 0
$ws.Cells.Item(174, 4).Value = 0
$ws.Cells.Item(174, 5).Value = 6
$ws.Cells.Item(174, 6).Value = 0
$ws.Cells.Item(174, 7).Value = 0
$ws.Cells.Item(174, 8).Value = 1

$ws.Cells.Item(175, 1).Value = 'Zimbabue'
$ws.Cells.Item(175, 2).Value = 7
$ws.Cells.Item(175, 3).Value = 0
$ws.Cells.Item(175, 4).Value = 0
$ws.Cells.Item(175, 5).Value = 6
$ws.Cells.Item(175, 6).Value = 0
$ws.Cells.Item(175, 7).Value = 0
$ws.Cells.Item(175, 8).Value = 1

$ws.Cells.Item(176, 1).Value = 'Santa Sede'
$ws.Cells.Item(176, 2).Value = 6
$ws.Cells.Item(176, 3).Value = 0
$ws.Cells.Item(176, 4).Value = 0
$ws.Cells.Item(176, 5).Value = 6
$ws.Cells.Item(176, 6).Value = 0
$ws.Cells.Item(176, 7).Value = 0
$ws.Cells.Item(176, 8).Value = 0

$ws.Cells.Item(177, 1).Value = 'San Martin (Parte Holandesa)'
$ws.Cells.Item(177, 2).Value = 6
$ws.Cells.Item(177, 3).Value = 0
$ws.Cells.Item(177, 4).Value = 0
$ws.Cells.Item(177, 5).Value = 6
$ws.Cells.Item(177, 6).Value = 0
$ws.Cells.Item(177, 7).Value = 0
$ws.Cells.Item(177, 8).Value = 0

$ws.Cells.Item(179, 1).Value = 'Eritrea'
$ws.Cells.Item(179, 2).Value = 6
$ws.Cells.Item(179, 3).Value = 0
$ws.Cells.Item(179, 4).Value = 0
$ws.Cells.Item(179, 5).Value = 6
$ws.Cells.Item(179, 6).Value = 0
$ws.Cells.Item(179, 7).Value = 0
$ws.Cells.Item(179, 8).Value = 0

$ws.Cells.Item(180, 1).Value = 'Siria'
$ws.Cells.Item(180, 2).Value = 6
$ws.Cells.Item(180, 3).Value = 1
$ws.Cells.Item(180, 4).Value = 0
$ws.Cells.Item(180, 5).Value = 5
$ws.Cells.Item(180, 6).Value = 0
$ws.Cells.Item(180, 7).Value = 1
$ws.Cells.Item(180, 8).Value = 1

$ws.Cells.Item(181, 1).Value = 'Cabo Verde'
$ws.Cells.Item(181, 2).Value = 6
$ws.Cells.Item(181, 3).Value = 0
$ws.Cells.Item(181, 4).Value = 0
$ws.Cells.Item(181, 5).Value = 5
$ws.Cells.Item(181, 6).Value = 0
$ws.Cells.Item(181, 7).Value = 0
$ws.Cells.Item(181, 8).Value = 1

$ws.Cells.Item(182, 1).Value = 'Montserrat'
$ws.Cells.Item(182, 2).Value = 5
$ws.Cells.Item(182, 3).Value = 0
$ws.Cells.Item(182, 4).Value = 0
$ws.Cells.Item(182, 5).Value = 5
$ws.Cells.Item(182, 6).Value = 0
$ws.Cells.Item(182, 7).Value = 0
$ws.Cells.Item(182, 8).Value = 0

$ws.Cells.Item(183, 1).Value = 'Angola'
$ws.Cells.Item(183, 2).Value = 5
$ws.Cells.Item(183, 3).Value = 0
$ws.Cells.Item(183, 4).Value = 0
$ws.Cells.Item(183, 5).Value = 5
$ws.Cells.Item(183, 6).Value = 0
$ws.Cells.Item(183, 7).Value = 0
$ws.Cells.Item(183, 8).Value = 0

$ws.Cells.Item(184, 1).Value = 'Mauritania'
$ws.Cells.Item(184, 2).Value = 5
$ws.Cells.Item(184, 3).Value = 0
$ws.Cells.Item(184, 4).Value = 0
$ws.Cells.Item(184, 5).Value = 5
$ws.Cells.Item(184, 6).Value = 0
$ws.Cells.Item(184, 7).Value = 0
$ws.Cells.Item(184, 8).Value = 0

$ws.Cells.Item(186, 1).Value = 'Fiyi'
$ws.Cells.Item(186, 2).Value = 5
$ws.Cells.Item(186, 3).Value = 0
$ws.Cells.Item(186, 4).Value = 0
$ws.Cells.Item(186, 5).Value = 5
$ws.Cells.Item(186, 6).Value = 0
$ws.Cells.Item(186, 7).Value = 0
$ws.Cells.Item(186, 8).Value = 0

$ws.Cells.Item(187, 1).Value = 'Nepal'
$ws.Cells.Item(187, 2).Value = 5
$ws.Cells.Item(187, 3).Value = 0
$ws.Cells.Item(187, 4).Value = 1
$ws.Cells.Item(187, 5).Value = 4
$ws.Cells.Item(187, 6).Value = 0
$ws.Cells.Item(187, 7).Value = 0
$ws.Cells.Item(187, 8).Value = 0

$ws.Cells.Item(188, 1).Value = 'Sudan'
$ws.Cells.Item(188, 2).Value = 5
$ws.Cells.Item(188, 3).Value = 0
$ws.Cells.Item(188, 4).Value = 0
$ws.Cells.Item(188, 5).Value = 4
$ws.Cells.Item(188, 6).Value = 0
$ws.Cells.Item(188, 7).Value = 0
$ws.Cells.Item(188, 8).Value = 1

$ws.Cells.Item(190, 1).Value = 'Islas Turcas y Caicos'
$ws.Cells.Item(190, 2).Value = 4
$ws.Cells.Item(190, 3).Value = 0
$ws.Cells.Item(190, 4).Value = 0
$ws.Cells.Item(190, 5).Value = 4
$ws.Cells.Item(190, 6).Value = 0
$ws.Cells.Item(190, 7).Value = 0
$ws.Cells.Item(190, 8).Value = 0

$ws.Cells.Item(191, 1).Value = 'Nicaragua'
$ws.Cells.Item(191, 2).Value = 4
$ws.Cells.Item(191, 3).Value = 0
$ws.Cells.Item(191, 4).Value = 0
$ws.Cells.Item(191, 5).Value = 3
$ws.Cells.Item(191, 6).Value = 0
$ws.Cells.Item(191, 7).Value = 0
$ws.Cells.Item(191, 8).Value = 1

$ws.Cells.Item(192, 1).Value = 'Santa Lucia'
$ws.Cells.Item(192, 2).Value = 3
$ws.Cells.Item(192, 3).Value = 1
$ws.Cells.Item(192, 4).Value = 1
$ws.Cells.Item(192, 5).Value = 3
$ws.Cells.Item(192, 6).Value = 0
$ws.Cells.Item(192, 7).Value = 0
$ws.Cells.Item(192, 8).Value = 0

$ws.Cells.Item(193, 1).Value = 'Somalia'
$ws.Cells.Item(193, 2).Value = 3
$ws.Cells.Item(193, 3).Value = 0
$ws.Cells.Item(193, 4).Value = 0
$ws.Cells.Item(193, 5).Value = 3
$ws.Cells.Item(193, 6).Value = 0
$ws.Cells.Item(193, 7).Value = 0
$ws.Cells.Item(193, 8).Value = 0

$ws.Cells.Item(194, 1).Value = 'Liberia'
$ws.Cells.Item(194, 2).Value = 3
$ws.Cells.Item(194, 3).Value = 0
$ws.Cells.Item(194, 4).Value = 0
$ws.Cells.Item(194, 5).Value = 3
$ws.Cells.Item(194, 6).Value = 0
$ws.Cells.Item(194, 7).Value = 0
$ws.Cells.Item(194, 8).Value = 0
